$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.65"
$ws.Range("E2").Value = "'-0.40%"
$ws.Range("E3").Value = "'-3.47%"
$ws.Range("D4").Value = "'5.291"
$ws.Range("E4").Value = "'1.94%"
$ws.Range("E5").Value = "'-0.40%"
$ws.Range("D6").Value = "'6.626"
$ws.Range("E6").Value = "'1.01%"
$ws.Range("D7").Value = "'3.211"
$ws.Range("E7").Value = "'3.45%"
$ws.Range("D8").Value = "'0.8544"
$ws.Range("E8").Value = "'-0.50%"
$ws.Range("D9").Value = "'0.8884"
$ws.Range("E9").Value = "'2.68%"
$ws.Range("D10").Value = "'0.1389"
$ws.Range("E10").Value = "'1.63%"
$ws.Range("D11").Value = "'0.07103"
$ws.Range("E11").Value = "'0.22%"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03321"
$ws.Range("E12").Value = "'-0.02%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03157"
$ws.Range("E13").Value = "'4.71%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09220"
$ws.Range("E14").Value = "'-1.77%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001528"
$ws.Range("E15").Value = "'-0.55%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0005988"
$ws.Range("E16").Value = "'-0.80%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005957"
$ws.Range("E17").Value = "'-1.77%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.499"
$ws.Range("E18").Value = "'0.21%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.172"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3166"
$ws.Range("E20").Value = "'-1.04%"
$ws.Range("D21").Value = "'0.1309"
$ws.Range("E21").Value = "'1.39%"
$ws.Range("D22").Value = "'3.499"
$ws.Range("E22").Value = "'0.48%"
$ws.Range("D23").Value = "'0.04071"
$ws.Range("E23").Value = "'-1.82%"
$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "'-0.21%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-0.39%"
$ws.Range("E26").Value = "'-16.73%"
$ws.Range("E27").Value = "'-0.94%"
$ws.Range("D40").Value = "'0.03793"
$ws.Range("E40").Value = "'1.01%"
$ws.Range("D41").Value = "'0.1067"
$ws.Range("E41").Value = "'-0.36%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002459"
$ws.Range("E42").Value = "'17.02%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.002949"
$ws.Range("E43").Value = "'-49.19%"
$ws.Range("D44").Value = "'0.009452"
$ws.Range("E44").Value = "'-0.18%"
$ws.Range("D45").Value = "'0.00005274"
$ws.Range("E45").Value = "'-0.34%"
$ws.Range("E46").Value = "'-0.12%"
$ws.Range("D47").Value = "'0.08906"
$ws.Range("E47").Value = "'56.09%"
$ws.Range("D48").Value = "'0.002257"
$ws.Range("E48").Value = "'-1.19%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.12%"
